$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "visit_type" column (column F) entirely; cells to the right
# shift one column to the left.
$ws.Range("F1:F1").EntireColumn.Delete()

# Update the active selection to match the new layout.
$ws.Range("F1").Select()
